# Applies the textual corrections described by the commit:
# "Small textual corrections to manual"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new, empty Title-styled paragraph before the existing
#    title paragraph ("Nuclear shape change analysis manual").
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$insertionPoint = $titlePara.Range.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Title'/><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr></w:p>")

# ---------------------------------------------------------------------
# 2. "... is the timed nuclear-channel images." -> "... is the
#    nuclear-channel image sequence."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The first image stack is the timed nuclear-channel images. This can be the original image or a pre-",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The first image stack is the nuclear-channel image sequence. This can be the original image or a pre-",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3. "The second required image stack is of labelled tracks." -> "The
#    second required image sequence contains the labelled tracks."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The second required image stack is of labelled tracks. Each instance of a specific nucleus needs",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The second required image sequence contains the labelled tracks. Each instance of a specific nucleus needs",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4. "(see Use below)" -> "(see 'Use' section below)" and "standard
#    name settings" -> "standard TITLE settings"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "the user has to manually identify the images (see Use below). If set to true, the macro will automatically identify the images by means of the standard name settings below.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the user has to manually identify the images (see " + [char]0x2018 + "Use" + [char]0x2019 + " section below). If set to true, the macro will automatically identify the images by means of the standard TITLE settings below.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 5. "For none binary nucleus image stacks" -> "For none-binary
#    nucleus image stacks"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "double thresholding. For none binary nucleus image stacks, leave this this setting on true.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "double thresholding. For none-binary nucleus image stacks, leave this this setting on true.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 6. "the NII score for this instance of the nucleus" -> "the NII
#    score for a particular nucleus"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The features not only contain the NII score for this instance of the nucleus, but also the Delta NII score",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The features not only contain the NII score for a particular nucleus, but also the Delta NII score",
    2) | Out-Null

Write-Output "edit complete"
